$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row (27) below the last existing data row (26), mirroring the
# pattern used by the other rows: VIN in column A, kit part number in column B.
$ws.Range("A27").Value = "4V4NC9EJ2EN168028"
$ws.Range("B27").Value = 100116

# Copy the style from the row above (B26) onto the new B27 cell so it keeps
# the same formatting as the rest of the kitPartNumber column.
$ws.Range("B26").Copy()
$ws.Range("B27").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Update the selection to match the new working range, as seen after
# selecting the newly added row together with the previous one.
$ws.Range("B26:B27").Select()
